# LoadConfigTestConfig.xlsx — add the "ExcelFiles" sheet (support for
# dispatching Excel files in config, alongside the existing TextFiles sheet)
# and move the active selection onto the new sheet.

$wb = $excel.ActiveWorkbook

# The previously-active sheet (TextFiles) will stop being the selected tab
# once ExcelFiles is added/activated; its own selection also moves from the
# old "next row" (A4) to the data block (A2:E3).
$tf = $wb.Worksheets.Item("TextFiles")
$tf.Range("A2:E3").Select()

# Insert the new sheet right after "TextFiles" so it becomes sheet index 5 /
# the last tab, and activate it (mirrors Excel's behavior of activating a
# newly inserted sheet).
$ws = $wb.Worksheets.Add($null, $tf)
$ws.Name = "ExcelFiles"

# Header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Path"
$ws.Range("C1").Value = "Folder"
$ws.Range("D1").Value = "Bucket"
$ws.Range("E1").Value = "Description"

# Row 2: local-file test case.
$ws.Range("A2").Value = "LocalFile"
$ws.Range("B2").Value = "Tests\Utility\LoadConfig\TestExcelFile.xlsx"
$ws.Range("E2").Value = "Testing a local file that should be loaded."

# Row 3: storage-bucket test case.
$ws.Range("A3").Value = "BucketFile"
$ws.Range("B3").Value = "TestExcelFile.xlsx"
$ws.Range("C3").Value = "LazyFramework"
$ws.Range("D3").Value = "LoadConfigTest"
$ws.Range("E3").Value = "Testing a storage bucket file that should be loaded."

# Match the saved selection on the new tab.
$ws.Range("G5").Select()
